$d = $word.ActiveDocument

# The second paragraph holds two M2Doc "complex" fields (fldChar begin /
# instrText* / fldChar end) that must become plain literal text of the
# form "{...}" made up of the same textual fragments, each kept in its
# own run (mirroring the original instrText run boundaries).
#
# Strategy: delete the two complex fields outright (Field.Delete removes
# the begin/instrText/end runs as a unit), which leaves plain Range
# addressing valid again, then re-insert the literal fragments as
# separate runs via repeated Range.InsertBefore calls on a collapsed
# Range (each call mints a new run; calling them back-to-front on the
# same collapsed position reproduces the original left-to-right order).

# Field 1 : {m:for | self.eClassifiers}
$field1Parts = @("{m", ":for ", "|", " ", "self.e", "Classifiers", "}")

# Field 2 : {m:endfor}
$field2Parts = @("{", "m", ":endfor}")

# Delete the first complex field (begin .. "m" ":for " "|" " " "self.e"
# "Classifiers" " " .. end).
$f1 = $d.Fields.Item(1)
$f1.Delete()

# The former second field is now the first (only) field left.
$f2 = $d.Fields.Item(1)
$f2.Delete()

$p2 = $d.Paragraphs.Item(2)
$pStart = $p2.Range.Start
$pEnd = $p2.Range.End

# Insert the "{m:endfor}" literal right before the paragraph mark (i.e.
# immediately after the "," that used to precede the second field).
for ($i = $field2Parts.Length - 1; $i -ge 0; $i--) {
    $r = $d.Range($pEnd - 1, $pEnd - 1)
    $r.InsertBefore($field2Parts[$i])
}

# Insert the "{m:for | self.eClassifiers}" literal at the very start of
# the paragraph (where the first field used to begin).
for ($i = $field1Parts.Length - 1; $i -ge 0; $i--) {
    $r = $d.Range($pStart, $pStart)
    $r.InsertBefore($field1Parts[$i])
}
